$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A50:B50").Copy() | Out-Null
$ws.Range("A52:B52").PasteSpecial(-4122) | Out-Null
$ws.Range("A50:B50").Copy() | Out-Null
$ws.Range("A53:B53").PasteSpecial(-4122) | Out-Null

$ws.Range("A52").Value = 45199
$ws.Range("B52").Value = "new times for embeddings"
$ws.Range("A53").Value = 45200
$ws.Range("B53").Value = "BA: doc2vec impl, doc2vec default model"

$ws.Rows.Item(52).RowHeight = 18
$ws.Rows.Item(53).RowHeight = 18

$ws.Range("C50").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 45
$excel.ActiveWindow.ScrollColumn = 1
